# Colocando header nos gráficos
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the same A-column layout:
#   A1 (new header), A2 "Hidro" ... A12 "GD"
# Add a header label in A1 (copying the bold/bordered style from B1),
# fix accented labels in A3/A4/A6/A8/A11, and drop the s="1" style that
# used to sit on A2:A12 (copying the unstyled format from B2).
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell A1, styled like the other header cells on row 1
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accented labels
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Strip the bold/bordered style from A2:A12 (now unstyled, like B2)
    $ws.Range("B2").Copy()
    $ws.Range("A2:A12").PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet 5 "Emissoes Totais (MtCO2eq)": add header, fix labels, drop the
# last row ("Teto") entirely.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("B2").Copy()
$ws5.Range("A2:A3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws5.Rows("4:4").Delete()

# ---------------------------------------------------------------------
# Sheet 6 "Custo Total (bilhões de R$)": add header, rename B1, update
# labels/values.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 used to hold the text "Custo"; it must still hold text ("2015"), not
# a number, so format it as text before assigning, then restore the
# original header style (copied from a cell that still carries it).
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 559

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

$ws6.Range("B2").Copy()
$ws6.Range("A2:A3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
